# Apply cryptocurrency price/volume updates and one row re-ranking
# (diff between before.xlsx and after.xlsx canonical OOXML)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.545.37"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "1.881.61"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.02"
$ws.Range("E4").Value = "  +0.55%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.43"
$ws.Range("E5").Value = "  +4.62%  "
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.15"
$ws.Range("E8").Value = "  +5.75%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("E10").Value = "  +1.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0994"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "2.150.81"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.34"
$ws.Range("E13").Value = "  +7.98%  "
$ws.Range("D14").Value = "1.920.49"
$ws.Range("E14").Value = "  +3.85%  "
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.79"
$ws.Range("E16").Value = "  +2.66%  "
$ws.Range("D17").Value = "35.581.55"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.46"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").Value = "0.0₃0808"
$ws.Range("E19").Value = "  +2.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "243.86"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  +1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.85"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("E24").Value = "  +1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.48"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.05"
$ws.Range("E26").Value = "  +25.96%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.32"
$ws.Range("E27").Value = "  +6.06%  "
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.974"
$ws.Range("E30").Value = "  +29.87%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0565"
$ws.Range("E31").Value = "  +2.04%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.08"
$ws.Range("E32").Value = "  +3.22%  "
$ws.Range("B33").Value = "BinanceUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.01"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.11"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.79"
$ws.Range("E35").Value = "  +11.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").Value = "  +5.53%  "
$ws.Range("E37").Value = "  +11.32%  "
$ws.Range("E38").Value = "  +2.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0205"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.25"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("D41").Value = "1.358.14"
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.22"
$ws.Range("E42").Value = "  +3.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0594"
$ws.Range("E43").Value = "  +11.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.34"
$ws.Range("E44").Value = "  +3.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.61"
$ws.Range("E45").Value = "  +43.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.42"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.74"
$ws.Range("E47").Value = "  +6.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.12"
$ws.Range("E49").Value = "  +33.26%  "
$ws.Range("D50").Value = "2.071.11"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0688"
$ws.Range("E51").Value = "  +2.57%  "
